# Saldo.xlsx update — "Add files via upload"
#
# Net effect (per the target diff):
#   - DIEGO's balance (account 004479965) is updated from 16173.58 to
#     32314.16 and the row is moved up, now appearing right after CEZAR
#     (004482090) and right before E3 (004267976).
#   - Five rows are dropped entirely: MARIANA (005000460), PEDRO
#     (005232019), CARLA (004643153), OTAVIO (004452946) and BHRUNA
#     (005295509).
#
# We operate on the original (pre-edit) row numbers, working from the
# bottom of the sheet upward so that each deletion/insertion is performed
# while the row numbers above it are still valid:
#   row 18     -> BHRUNA                              (delete)
#   rows 9-12  -> MARIANA, PEDRO, CARLA, OTAVIO        (delete block)
#   row 7      -> old DIEGO row (16173.58)             (delete)
#   row 5      -> insert new DIEGO row (32314.16) just above E3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove BHRUNA (005295509 / 569.57)
$ws.Rows.Item(18).Delete()

# 2) Remove the MARIANA / PEDRO / CARLA / OTAVIO block
$ws.Range("9:12").Delete()

# 3) Remove the old DIEGO row (account 004479965, balance 16173.58)
$ws.Rows.Item(7).Delete()

# 4) Insert the updated DIEGO row right before E3 (004267976), with the
#    new balance. The leading apostrophe forces the account number to be
#    stored as text (preserving the leading zeros), matching the rest of
#    column A.
$ws.Rows.Item(5).Insert()
$ws.Cells.Item(5, 1).Value = "'004479965"
$ws.Cells.Item(5, 2).Value = "DIEGO"
$ws.Cells.Item(5, 3).Value = 32314.16
